# Simple refactor of delete: add a small "R:null" label textbox to the
# "Delete Node in Sorted tree" slide (the 4th slide, sldId 261), matching
# the caption style already used by the other textboxes on that slide.
#
# We duplicate an existing plain (non-bold) caption textbox on the same
# slide ("TextBox 6", which holds "D: the node to be deleted / P: the
# parent node of the node to be deleted") so the new shape inherits the
# same body/paragraph/run formatting already baked into the slide, then
# reposition, resize, rename and retext it to become the new label.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

$src = $s.Shapes.Item(5)   # "TextBox 6" - D:/P: caption textbox
$dupRange = $src.Duplicate()
$shp = $dupRange.Item(1)

$shp.Name = "TextBox 1"

# Target position/size are given in EMU in the canonical OOXML; the
# PowerPoint COM object model works in points (1 pt = 12700 EMU), so we
# convert by dividing. A tiny epsilon nudge on Width compensates for
# floating point round-trip rounding so the resulting EMU value lands on
# the exact target (763351) rather than one unit short.
$shp.Left   = 1426029.0 / 12700.0
$shp.Top    = 5965371.0 / 12700.0
$shp.Width  = 763351.0  / 12700.0 + 0.00001
$shp.Height = 369332.0  / 12700.0

$shp.TextFrame.TextRange.Text = "R:null"
